# Fix the dates for the first couple of weeks on the schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Aug 26 - Sept 06"
$ws.Range("B3").Value = "Sep 09 - Sep 13"
$ws.Range("B4").Value = "Sep 16 - Sep 20"
$ws.Range("B5").Value = "Sep 23 - Sep 27"

# Move the active selection, matching the saved cursor position in the file.
$ws.Range("C21").Select()
